# "Getting ready to send out to collaborators"
#
# 1. Update the dataset Description on the Collection sheet to mention the
#    extra directory of text files.
# 2. Add a new RELATION:Keywords row to the Collection sheet.
# 3. Make the Collection sheet the active/selected tab (instead of Places),
#    with B4 as the selected cell, matching the refreshed view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection")

# --- Content edits -------------------------------------------------------

$ws.Range("B3").Value = "This is a simple dataset for demonstration purposes it contains just one image and a directory full of useless text files."

$ws.Range("A12").Value = "RELATION:Keywords"
$ws.Range("B12").Value = "Dogs, Fences, The Gully"

# --- View state: make Collection the active tab/selection ----------------

$ws.Activate()
$ws.Range("B4").Select() | Out-Null

$wb.Windows.Item(1).TabRatio = 955 | Out-Null
